# Applies "Update countries & provincias Spain" edit:
#  - refresh the "Datos actualizados" timestamp string
#  - fix 5 mis-ordered country names (A column) that had been swapped
#  - update the Covid-19 stat columns (B-H) with refreshed figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Refresh "last updated" timestamp (row 1, col A)
$ws.Cells.Item(1, 1).Value2 = "Datos actualizados a 18 de Julio de 2020 a las 20:26"

# 2) Correct the country names that were out of order in column A
$countryNameFixes = @{
    44 = "Israel"
    45 = "Portugal"
    46 = "Singapur"
    79 = "Etiopia"
    80 = "Noruega"
    81 = "Republica de Macedonia"
    112 = "Libano"
    113 = "Sri Lanka"
    121 = "Guinea-Bisau"
    122 = "Eslovenia"
    123 = "Cabo Verde"
    144 = "Republica de Chipre"
    145 = "Uruguay"
}
foreach ($row in $countryNameFixes.Keys) {
    $ws.Cells.Item([int]$row, 1).Value2 = $countryNameFixes[$row]
}

# 3) Update the statistic columns (B=Casos totales, C=Nuevos casos,
#    D=Casos activos, E=Recuperados, G=Muertes hoy, H=Muertes) with new figures
$cellUpdates = @(
    @(4, 2, 3798407),
    @(4, 3, 28395),
    @(4, 4, 1756504),
    @(4, 5, 1899431),
    @(4, 7, 408),
    @(4, 8, 142472),
    @(6, 2, 1076535),
    @(6, 3, 36078),
    @(6, 4, 677480),
    @(6, 5, 372229),
    @(6, 7, 541),
    @(6, 8, 26826),
    @(18, 2, 218717),
    @(18, 3, 918),
    @(18, 4, 201013),
    @(18, 5, 12229),
    @(18, 7, 17),
    @(18, 8, 5475),
    @(44, 2, 49204),
    @(44, 3, 1745),
    @(44, 4, 21302),
    @(44, 5, 27502),
    @(44, 7, 8),
    @(44, 8, 400),
    @(45, 2, 48390),
    @(45, 3, 313),
    @(45, 4, 33153),
    @(45, 5, 13553),
    @(45, 7, 2),
    @(45, 8, 1684),
    @(46, 2, 47655),
    @(46, 3, 202),
    @(46, 4, 43833),
    @(46, 5, 3795),
    @(46, 8, 27),
    @(51, 2, 35301),
    @(51, 3, 72),
    @(51, 4, 23273),
    @(51, 5, 10864),
    @(51, 7, 17),
    @(51, 8, 1164),
    @(58, 2, 25750),
    @(58, 3, 20),
    @(58, 5, 633),
    @(58, 7, 1),
    @(58, 8, 1753),
    @(60, 4, 15744),
    @(60, 5, 5737),
    @(65, 2, 17015),
    @(65, 3, 289),
    @(65, 4, 14620),
    @(65, 5, 2126),
    @(65, 7, 5),
    @(65, 8, 269),
    @(79, 2, 9147),
    @(79, 3, 344),
    @(79, 4, 2430),
    @(79, 5, 6554),
    @(79, 7, 13),
    @(79, 8, 163),
    @(80, 2, 9028),
    @(80, 3, 3),
    @(80, 4, 8138),
    @(80, 5, 635),
    @(80, 7, 0),
    @(80, 8, 255),
    @(81, 2, 9026),
    @(81, 3, 240),
    @(81, 4, 4727),
    @(81, 5, 3885),
    @(81, 7, 8),
    @(81, 8, 414),
    @(86, 2, 8204),
    @(86, 3, 440),
    @(86, 7, 6),
    @(86, 8, 59),
    @(92, 2, 6655),
    @(92, 3, 146),
    @(92, 4, 4272),
    @(92, 5, 2346),
    @(92, 7, 3),
    @(92, 8, 37),
    @(109, 2, 2930),
    @(109, 3, 17),
    @(109, 4, 2354),
    @(109, 5, 561),
    @(112, 2, 2775),
    @(112, 3, 75),
    @(112, 4, 1485),
    @(112, 5, 1250),
    @(112, 8, 40),
    @(113, 2, 2703),
    @(113, 3, 6),
    @(113, 4, 2023),
    @(113, 5, 669),
    @(113, 8, 11),
    @(116, 2, 2445),
    @(116, 3, 1),
    @(116, 4, 2304),
    @(116, 5, 54),
    @(121, 2, 1949),
    @(121, 3, 22),
    @(121, 4, 803),
    @(121, 5, 1120),
    @(121, 8, 26),
    @(122, 2, 1940),
    @(122, 3, 24),
    @(122, 4, 1568),
    @(122, 5, 261),
    @(122, 8, 111),
    @(123, 2, 1939),
    @(123, 4, 902),
    @(123, 5, 1018),
    @(123, 8, 19),
    @(144, 3, 4),
    @(144, 4, 845),
    @(144, 5, 173),
    @(144, 8, 19),
    @(145, 2, 1037),
    @(145, 4, 917),
    @(145, 5, 88),
    @(145, 8, 32),
    @(157, 2, 624),
    @(157, 3, 10),
    @(157, 5, 149),
    @(165, 2, 340),
    @(165, 3, 1),
    @(165, 4, 273),
    @(165, 5, 61)
)
foreach ($u in $cellUpdates) {
    $ws.Cells.Item($u[0], $u[1]).Value2 = $u[2]
}
